# Auto-generated edit script applying the diff to before.xlsx
$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(19, 8).Value = 2113.8333
$ws.Cells.Item(19, 9).Value = 2227.6667
$ws.Cells.Item(19, 10).Value = 2000
$ws.Cells.Item(19, 11).Value = 2227.6667
$ws.Cells.Item(19, 12).Value = 2000
$ws.Cells.Item(19, 13).Value = -2052.6667
$ws.Cells.Item(19, 14).Value = -2350
$ws.Cells.Item(32, 8).Value = 3514.1538
$ws.Cells.Item(32, 10).Value = 4068.5
$ws.Cells.Item(32, 12).Value = 4068.5
$ws.Cells.Item(32, 14).Value = -4720.5
$ws.Cells.Item(43, 8).Value = 835491.5
$ws.Cells.Item(43, 9).Value = 2244.2222
$ws.Cells.Item(43, 10).Value = 3335233.2
$ws.Cells.Item(43, 11).Value = 2244.2222
$ws.Cells.Item(43, 12).Value = 3335233.2
$ws.Cells.Item(43, 13).Value = -2175.2222
$ws.Cells.Item(43, 14).Value = -3335371.2
$ws.Cells.Item(51, 8).Value = 9552.111000000001
$ws.Cells.Item(51, 9).Value = 4995
$ws.Cells.Item(51, 10).Value = 15248.5
$ws.Cells.Item(51, 11).Value = 4995
$ws.Cells.Item(51, 12).Value = 15248.5
$ws.Cells.Item(51, 13).Value = -4511
$ws.Cells.Item(51, 14).Value = -16216.5
$ws.Cells.Item(62, 8).Value = 1665.3334
$ws.Cells.Item(62, 9).Value = 998.5
$ws.Cells.Item(62, 10).Value = 2999
$ws.Cells.Item(62, 11).Value = 998.5
$ws.Cells.Item(62, 12).Value = 2999
$ws.Cells.Item(62, 13).Value = -374.5
$ws.Cells.Item(62, 14).Value = -4247
$ws.Cells.Item(64, 8).Value = 6587.3
$ws.Cells.Item(64, 10).Value = 6962.6665
$ws.Cells.Item(64, 12).Value = 6962.6665
$ws.Cells.Item(64, 14).Value = -7458.6665
$ws.Cells.Item(65, 8).Value = 1665.3334
$ws.Cells.Item(65, 9).Value = 998.5
$ws.Cells.Item(65, 10).Value = 2999
$ws.Cells.Item(65, 11).Value = 4992.5
$ws.Cells.Item(65, 12).Value = 14995
$ws.Cells.Item(65, 13).Value = -1872.5
$ws.Cells.Item(65, 14).Value = -21235
$ws.Cells.Item(67, 8).Value = 6587.3
$ws.Cells.Item(67, 10).Value = 6962.6665
$ws.Cells.Item(67, 12).Value = 6962.6665
$ws.Cells.Item(67, 14).Value = -8678.666499999999
$ws.Cells.Item(74, 8).Value = 5695.1333
$ws.Cells.Item(74, 9).Value = 4697.6665
$ws.Cells.Item(74, 10).Value = 5944.5
$ws.Cells.Item(74, 11).Value = 4697.6665
$ws.Cells.Item(74, 12).Value = 5944.5
$ws.Cells.Item(74, 13).Value = -3761.6665
$ws.Cells.Item(74, 14).Value = -7816.5
$ws.Cells.Item(77, 8).Value = 5695.1333
$ws.Cells.Item(77, 9).Value = 4697.6665
$ws.Cells.Item(77, 10).Value = 5944.5
$ws.Cells.Item(77, 11).Value = 23488.3325
$ws.Cells.Item(77, 12).Value = 29722.5
$ws.Cells.Item(77, 13).Value = -18808.3325
$ws.Cells.Item(77, 14).Value = -39082.5
$ws.Cells.Item(92, 8).Value = 802.2727
$ws.Cells.Item(92, 9).Value = 802.2727
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 802.2727
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 445.7273
$ws.Cells.Item(92, 14).ClearContents()
$ws.Cells.Item(112, 8).Value = 2511.5
$ws.Cells.Item(112, 9).Value = 1893.8334
$ws.Cells.Item(112, 10).Value = 2974.75
$ws.Cells.Item(112, 11).Value = 5681.5002
$ws.Cells.Item(112, 12).Value = 8924.25
$ws.Cells.Item(112, 13).Value = -4573.5002
$ws.Cells.Item(112, 14).Value = -11140.25
$ws.Cells.Item(116, 8).Value = 7359.375
$ws.Cells.Item(116, 9).Value = 7359.375
$ws.Cells.Item(116, 11).Value = 7359.375
$ws.Cells.Item(116, 13).Value = -3917.375
$ws.Cells.Item(123, 8).Value = 58100
$ws.Cells.Item(123, 10).Value = 58100
$ws.Cells.Item(123, 12).Value = 58100
$ws.Cells.Item(123, 14).Value = -67900
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 13).ClearContents()
$ws.Cells.Item(125, 8).Value = 5331
$ws.Cells.Item(125, 9).Value = 5719.8
$ws.Cells.Item(125, 11).Value = 51478.2
$ws.Cells.Item(125, 13).Value = -49018.2
$ws.Cells.Item(132, 8).Value = 4865
$ws.Cells.Item(132, 9).Value = 4852.2383
$ws.Cells.Item(132, 10).Value = 4999
$ws.Cells.Item(132, 11).Value = 14556.7149
$ws.Cells.Item(132, 12).Value = 14997
$ws.Cells.Item(132, 13).Value = -12026.7149
$ws.Cells.Item(132, 14).Value = -20057
$ws.Cells.Item(137, 8).Value = 16392.777
$ws.Cells.Item(137, 9).Value = 23079.732
$ws.Cells.Item(137, 10).Value = 8034.0835
$ws.Cells.Item(137, 11).Value = 69239.196
$ws.Cells.Item(137, 12).Value = 24102.2505
$ws.Cells.Item(137, 13).Value = -66689.196
$ws.Cells.Item(137, 14).Value = -29202.2505
$ws.Cells.Item(138, 8).Value = 5189.8335
$ws.Cells.Item(138, 9).Value = 1184.0714
$ws.Cells.Item(138, 10).Value = 7738.9546
$ws.Cells.Item(138, 11).Value = 3552.2142
$ws.Cells.Item(138, 12).Value = 23216.8638
$ws.Cells.Item(138, 13).Value = 1587.7858
$ws.Cells.Item(138, 14).Value = -33496.8638
$ws.Cells.Item(141, 8).Value = 6928.4
$ws.Cells.Item(141, 9).Value = 7884.25
$ws.Cells.Item(141, 11).Value = 23652.75
$ws.Cells.Item(141, 13).Value = -18472.75

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(19, 8).Value = 845
$ws.Cells.Item(19, 9).Value = 490
$ws.Cells.Item(19, 10).Value = 1200
$ws.Cells.Item(19, 11).Value = 490
$ws.Cells.Item(19, 12).Value = 1200
$ws.Cells.Item(19, 13).Value = -261
$ws.Cells.Item(19, 14).Value = -1658
$ws.Cells.Item(32, 8).Value = 6269.375
$ws.Cells.Item(32, 9).Value = 5826.4517
$ws.Cells.Item(32, 11).Value = 5826.4517
$ws.Cells.Item(32, 13).Value = -5539.4517
$ws.Cells.Item(43, 8).Value = 500029980
$ws.Cells.Item(43, 10).Value = 59999
$ws.Cells.Item(43, 12).Value = 59999
$ws.Cells.Item(43, 14).Value = -60625
$ws.Cells.Item(45, 8).Value = 3255.7307
$ws.Cells.Item(45, 9).Value = 2954.3076
$ws.Cells.Item(45, 10).Value = 3557.1538
$ws.Cells.Item(45, 11).Value = 2954.3076
$ws.Cells.Item(45, 12).Value = 3557.1538
$ws.Cells.Item(45, 13).Value = -2577.3076
$ws.Cells.Item(45, 14).Value = -4311.1538
$ws.Cells.Item(61, 8).Value = 2336.318
$ws.Cells.Item(61, 9).Value = 2026.2632
$ws.Cells.Item(61, 10).Value = 4300
$ws.Cells.Item(61, 11).Value = 2026.2632
$ws.Cells.Item(61, 12).Value = 4300
$ws.Cells.Item(61, 13).Value = -1814.2632
$ws.Cells.Item(61, 14).Value = -4724
$ws.Cells.Item(63, 8).Value = 2164
$ws.Cells.Item(63, 9).Value = 2164
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 2164
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = -1478
$ws.Cells.Item(63, 14).ClearContents()
$ws.Cells.Item(66, 8).Value = 2164
$ws.Cells.Item(66, 9).Value = 2164
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 10820
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 13).Value = -7388
$ws.Cells.Item(66, 14).ClearContents()
$ws.Cells.Item(74, 8).Value = 25558.166
$ws.Cells.Item(74, 9).Value = 26084.299
$ws.Cells.Item(74, 11).Value = 26084.299
$ws.Cells.Item(74, 13).Value = -25210.299
$ws.Cells.Item(77, 8).Value = 25558.166
$ws.Cells.Item(77, 9).Value = 26084.299
$ws.Cells.Item(77, 11).Value = 130421.495
$ws.Cells.Item(77, 13).Value = -126053.495
$ws.Cells.Item(97, 8).Value = 1039.7368
$ws.Cells.Item(97, 9).Value = 907.0625
$ws.Cells.Item(97, 11).Value = 907.0625
$ws.Cells.Item(97, 13).Value = -411.0625
$ws.Cells.Item(110, 8).Value = 32512.455
$ws.Cells.Item(110, 9).Value = 35693.7
$ws.Cells.Item(110, 10).Value = 700
$ws.Cells.Item(110, 11).Value = 35693.7
$ws.Cells.Item(110, 12).Value = 700
$ws.Cells.Item(110, 13).Value = -33648.7
$ws.Cells.Item(110, 14).Value = -4790
$ws.Cells.Item(132, 8).Value = 75786.85000000001
$ws.Cells.Item(132, 9).Value = 6179.8335
$ws.Cells.Item(132, 10).Value = 702250
$ws.Cells.Item(132, 11).Value = 18539.5005
$ws.Cells.Item(132, 12).Value = 2106750
$ws.Cells.Item(132, 13).Value = -16009.5005
$ws.Cells.Item(132, 14).Value = -2111810
$ws.Cells.Item(136, 8).Value = 2336.318
$ws.Cells.Item(136, 9).Value = 2026.2632
$ws.Cells.Item(136, 10).Value = 4300
$ws.Cells.Item(136, 11).Value = 6078.7896
$ws.Cells.Item(136, 12).Value = 12900
$ws.Cells.Item(136, 13).Value = -3528.7896
$ws.Cells.Item(136, 14).Value = -18000

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(22, 8).Value = 837.25
$ws.Cells.Item(22, 9).Value = 699.6667
$ws.Cells.Item(22, 11).Value = 699.6667
$ws.Cells.Item(22, 13).Value = -526.6667
$ws.Cells.Item(86, 8).Value = 59462.08
$ws.Cells.Item(86, 9).Value = 49625
$ws.Cells.Item(86, 10).Value = 75201.39999999999
$ws.Cells.Item(86, 11).Value = 49625
$ws.Cells.Item(86, 12).Value = 75201.39999999999
$ws.Cells.Item(86, 13).Value = -48502
$ws.Cells.Item(86, 14).Value = -77447.39999999999
$ws.Cells.Item(89, 8).Value = 59462.08
$ws.Cells.Item(89, 9).Value = 49625
$ws.Cells.Item(89, 10).Value = 75201.39999999999
$ws.Cells.Item(89, 11).Value = 248125
$ws.Cells.Item(89, 12).Value = 376007
$ws.Cells.Item(89, 13).Value = -242509
$ws.Cells.Item(89, 14).Value = -387239
$ws.Cells.Item(96, 8).Value = 40412.4
$ws.Cells.Item(96, 9).Value = 7246.2856
$ws.Cells.Item(96, 11).Value = 7246.2856
$ws.Cells.Item(96, 13).Value = -4500.2856
$ws.Cells.Item(102, 8).Value = 64528
$ws.Cells.Item(102, 9).Value = 23778
$ws.Cells.Item(102, 10).Value = 84903
$ws.Cells.Item(102, 11).Value = 23778
$ws.Cells.Item(102, 12).Value = 84903
$ws.Cells.Item(102, 13).Value = -20533
$ws.Cells.Item(102, 14).Value = -91393
$ws.Cells.Item(134, 8).Value = 3745
$ws.Cells.Item(134, 9).Value = 3745
$ws.Cells.Item(134, 11).Value = 11235
$ws.Cells.Item(134, 13).Value = -8700

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 2218.375
$ws.Cells.Item(16, 10).Value = 2122
$ws.Cells.Item(16, 12).Value = 2122
$ws.Cells.Item(16, 14).Value = -2696
$ws.Cells.Item(31, 8).Value = 9262845
$ws.Cells.Item(31, 9).Value = 2727.2104
$ws.Cells.Item(31, 10).Value = 31255624
$ws.Cells.Item(31, 11).Value = 2727.2104
$ws.Cells.Item(31, 12).Value = 31255624
$ws.Cells.Item(31, 13).Value = -2432.2104
$ws.Cells.Item(31, 14).Value = -31256214
$ws.Cells.Item(34, 8).Value = 9262845
$ws.Cells.Item(34, 9).Value = 2727.2104
$ws.Cells.Item(34, 10).Value = 31255624
$ws.Cells.Item(34, 11).Value = 2727.2104
$ws.Cells.Item(34, 12).Value = 31255624
$ws.Cells.Item(34, 13).Value = -2525.2104
$ws.Cells.Item(34, 14).Value = -31256028
$ws.Cells.Item(58, 8).Value = 39290628
$ws.Cells.Item(58, 9).Value = 7700
$ws.Cells.Item(58, 11).Value = 7700
$ws.Cells.Item(58, 13).Value = -7497
$ws.Cells.Item(62, 8).Value = 3621.1875
$ws.Cells.Item(62, 9).Value = 3485.7144
$ws.Cells.Item(62, 11).Value = 3485.7144
$ws.Cells.Item(62, 13).Value = -2861.7144
$ws.Cells.Item(65, 8).Value = 3621.1875
$ws.Cells.Item(65, 9).Value = 3485.7144
$ws.Cells.Item(65, 11).Value = 17428.572
$ws.Cells.Item(65, 13).Value = -14308.572
$ws.Cells.Item(86, 8).Value = 3926.7144
$ws.Cells.Item(86, 9).Value = 3898
$ws.Cells.Item(86, 10).Value = 3998.5
$ws.Cells.Item(86, 11).Value = 3898
$ws.Cells.Item(86, 12).Value = 3998.5
$ws.Cells.Item(86, 13).Value = -2775
$ws.Cells.Item(86, 14).Value = -6244.5
$ws.Cells.Item(89, 8).Value = 3926.7144
$ws.Cells.Item(89, 9).Value = 3898
$ws.Cells.Item(89, 10).Value = 3998.5
$ws.Cells.Item(89, 11).Value = 19490
$ws.Cells.Item(89, 12).Value = 19992.5
$ws.Cells.Item(89, 13).Value = -13874
$ws.Cells.Item(89, 14).Value = -31224.5
$ws.Cells.Item(99, 8).Value = 6012.375
$ws.Cells.Item(99, 9).Value = 6299.857
$ws.Cells.Item(99, 11).Value = 6299.857
$ws.Cells.Item(99, 13).Value = -4801.857
$ws.Cells.Item(107, 8).Value = 656.4545000000001
$ws.Cells.Item(107, 9).Value = 634.25
$ws.Cells.Item(107, 11).Value = 634.25
$ws.Cells.Item(107, 13).Value = 1285.75
$ws.Cells.Item(113, 8).Value = 2218.375
$ws.Cells.Item(113, 10).Value = 2122
$ws.Cells.Item(113, 12).Value = 2122
$ws.Cells.Item(113, 14).Value = -6462
$ws.Cells.Item(122, 8).Value = 4050649.5
$ws.Cells.Item(122, 9).Value = 1611
$ws.Cells.Item(122, 10).Value = 10529111
$ws.Cells.Item(122, 11).Value = 4833
$ws.Cells.Item(122, 12).Value = 31587333
$ws.Cells.Item(122, 13).Value = -2383
$ws.Cells.Item(122, 14).Value = -31592233
$ws.Cells.Item(126, 8).Value = 6012.375
$ws.Cells.Item(126, 9).Value = 6299.857
$ws.Cells.Item(126, 11).Value = 18899.571
$ws.Cells.Item(126, 13).Value = -16429.571
$ws.Cells.Item(132, 8).Value = 3873.1482
$ws.Cells.Item(132, 9).Value = 3791.4
$ws.Cells.Item(132, 10).Value = 4895
$ws.Cells.Item(132, 11).Value = 11374.2
$ws.Cells.Item(132, 12).Value = 14685
$ws.Cells.Item(132, 13).Value = -8844.200000000001
$ws.Cells.Item(132, 14).Value = -19745
$ws.Cells.Item(134, 8).Value = 3853.182
$ws.Cells.Item(134, 9).Value = 3192
$ws.Cells.Item(134, 10).Value = 5616.3335
$ws.Cells.Item(134, 11).Value = 9576
$ws.Cells.Item(134, 12).Value = 16849.0005
$ws.Cells.Item(134, 13).Value = -7041
$ws.Cells.Item(134, 14).Value = -21919.0005
$ws.Cells.Item(136, 8).Value = 39290628
$ws.Cells.Item(136, 9).Value = 7700
$ws.Cells.Item(136, 11).Value = 23100
$ws.Cells.Item(136, 13).Value = -20550
$ws.Cells.Item(140, 8).Value = 124999.5
$ws.Cells.Item(140, 10).Value = 124999.5
$ws.Cells.Item(140, 12).Value = 124999.5
$ws.Cells.Item(140, 14).Value = -135359.5

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(4, 8).Value = 46207156
$ws.Cells.Item(4, 9).Value = 107284580
$ws.Cells.Item(4, 10).Value = 17122672
$ws.Cells.Item(4, 11).Value = 321853740
$ws.Cells.Item(4, 12).Value = 51368016
$ws.Cells.Item(4, 13).Value = -321853628
$ws.Cells.Item(4, 14).Value = -51368240
$ws.Cells.Item(5, 8).Value = 1827.6
$ws.Cells.Item(5, 10).Value = 4299.75
$ws.Cells.Item(5, 12).Value = 12899.25
$ws.Cells.Item(5, 14).Value = -13123.25
$ws.Cells.Item(12, 8).Value = 436.9565
$ws.Cells.Item(12, 9).Value = 197.28572
$ws.Cells.Item(12, 10).Value = 541.8125
$ws.Cells.Item(12, 11).Value = 591.85716
$ws.Cells.Item(12, 12).Value = 1625.4375
$ws.Cells.Item(12, 13).Value = -418.85716
$ws.Cells.Item(12, 14).Value = -1971.4375
$ws.Cells.Item(61, 8).Value = 137.4
$ws.Cells.Item(61, 9).Value = 144.38461
$ws.Cells.Item(61, 10).Value = 92
$ws.Cells.Item(61, 11).Value = 433.15383
$ws.Cells.Item(61, 12).Value = 276
$ws.Cells.Item(61, 13).Value = -218.15383
$ws.Cells.Item(61, 14).Value = -706
$ws.Cells.Item(92, 8).Value = 1148.6
$ws.Cells.Item(92, 9).Value = 1300
$ws.Cells.Item(92, 10).Value = 1110.75
$ws.Cells.Item(92, 11).Value = 3900
$ws.Cells.Item(92, 12).Value = 3332.25
$ws.Cells.Item(92, 13).Value = -2652
$ws.Cells.Item(92, 14).Value = -5828.25
$ws.Cells.Item(122, 8).Value = 1227.4706
$ws.Cells.Item(122, 10).Value = 1438.5834
$ws.Cells.Item(122, 12).Value = 12947.2506
$ws.Cells.Item(122, 14).Value = -17847.2506
$ws.Cells.Item(132, 8).Value = 5290.5
$ws.Cells.Item(132, 9).Value = 1264.8
$ws.Cells.Item(132, 11).Value = 11383.2
$ws.Cells.Item(132, 13).Value = -8853.199999999999
$ws.Cells.Item(135, 8).Value = 1827.6
$ws.Cells.Item(135, 10).Value = 4299.75
$ws.Cells.Item(135, 12).Value = 38697.75
$ws.Cells.Item(135, 14).Value = -43767.75
$ws.Cells.Item(139, 8).Value = 2782
$ws.Cells.Item(139, 9).Value = 2615.2
$ws.Cells.Item(139, 11).Value = 7845.599999999999
$ws.Cells.Item(139, 13).Value = -2705.599999999999

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 155.08696
$ws.Cells.Item(2, 9).Value = 60.866665
$ws.Cells.Item(2, 10).Value = 331.75
$ws.Cells.Item(2, 11).Value = 60.866665
$ws.Cells.Item(2, 12).Value = 331.75
$ws.Cells.Item(2, 13).Value = 52.133335
$ws.Cells.Item(2, 14).Value = -557.75
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 13).ClearContents()
$ws.Cells.Item(55, 8).Value = 30000
$ws.Cells.Item(55, 10).Value = 30000
$ws.Cells.Item(55, 12).Value = 30000
$ws.Cells.Item(55, 14).Value = -30654
$ws.Cells.Item(70, 8).Value = 6848.706
$ws.Cells.Item(70, 9).Value = 4751.4165
$ws.Cells.Item(70, 11).Value = 4751.4165
$ws.Cells.Item(70, 13).Value = -4481.4165
$ws.Cells.Item(73, 8).Value = 6848.706
$ws.Cells.Item(73, 9).Value = 4751.4165
$ws.Cells.Item(73, 11).Value = 4751.4165
$ws.Cells.Item(73, 13).Value = -3815.4165
$ws.Cells.Item(80, 8).Value = 6314.4165
$ws.Cells.Item(80, 10).Value = 7548.857
$ws.Cells.Item(80, 12).Value = 7548.857
$ws.Cells.Item(80, 14).Value = -9544.857
$ws.Cells.Item(83, 8).Value = 6314.4165
$ws.Cells.Item(83, 10).Value = 7548.857
$ws.Cells.Item(83, 12).Value = 37744.285
$ws.Cells.Item(83, 14).Value = -47728.285
$ws.Cells.Item(102, 8).Value = 5415.875
$ws.Cells.Item(102, 9).Value = 1974
$ws.Cells.Item(102, 10).Value = 20330.666
$ws.Cells.Item(102, 11).Value = 1974
$ws.Cells.Item(102, 12).Value = 20330.666
$ws.Cells.Item(102, 13).Value = -352
$ws.Cells.Item(102, 14).Value = -23574.666
$ws.Cells.Item(107, 8).Value = 896.8333
$ws.Cells.Item(107, 10).Value = 1579.8
$ws.Cells.Item(107, 12).Value = 1579.8
$ws.Cells.Item(107, 14).Value = -5419.8
$ws.Cells.Item(113, 8).Value = 2000
$ws.Cells.Item(113, 9).Value = 1000
$ws.Cells.Item(113, 10).Value = 2333.3333
$ws.Cells.Item(113, 11).Value = 1000
$ws.Cells.Item(113, 12).Value = 2333.3333
$ws.Cells.Item(113, 13).Value = 1170
$ws.Cells.Item(113, 14).Value = -6673.3333
$ws.Cells.Item(122, 8).Value = 16669306
$ws.Cells.Item(122, 9).Value = 1954
$ws.Cells.Item(122, 10).Value = 41670332
$ws.Cells.Item(122, 11).Value = 5862
$ws.Cells.Item(122, 12).Value = 125010996
$ws.Cells.Item(122, 13).Value = -3412
$ws.Cells.Item(122, 14).Value = -125015896
$ws.Cells.Item(126, 8).Value = 1435.4286
$ws.Cells.Item(126, 9).Value = 1536.3334
$ws.Cells.Item(126, 10).Value = 830
$ws.Cells.Item(126, 11).Value = 4609.0002
$ws.Cells.Item(126, 12).Value = 2490
$ws.Cells.Item(126, 13).Value = -2139.0002
$ws.Cells.Item(126, 14).Value = -7430
$ws.Cells.Item(132, 8).Value = 2903.5833
$ws.Cells.Item(132, 9).Value = 2076.6365
$ws.Cells.Item(132, 11).Value = 6229.9095
$ws.Cells.Item(132, 13).Value = -3699.9095

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 14).ClearContents()
$ws.Cells.Item(40, 8).Value = 4470
$ws.Cells.Item(40, 9).Value = 3593
$ws.Cells.Item(40, 10).Value = 6662.5
$ws.Cells.Item(40, 11).Value = 3593
$ws.Cells.Item(40, 12).Value = 6662.5
$ws.Cells.Item(40, 13).Value = -3457
$ws.Cells.Item(40, 14).Value = -6934.5
$ws.Cells.Item(45, 8).Value = 29917
$ws.Cells.Item(45, 9).Value = 14950
$ws.Cells.Item(45, 10).Value = 37400.5
$ws.Cells.Item(45, 11).Value = 14950
$ws.Cells.Item(45, 12).Value = 37400.5
$ws.Cells.Item(45, 13).Value = -14543
$ws.Cells.Item(45, 14).Value = -38214.5
$ws.Cells.Item(46, 8).Value = 2226.44
$ws.Cells.Item(46, 9).Value = 956.2857
$ws.Cells.Item(46, 10).Value = 3843
$ws.Cells.Item(46, 11).Value = 956.2857
$ws.Cells.Item(46, 12).Value = 3843
$ws.Cells.Item(46, 13).Value = -768.2857
$ws.Cells.Item(46, 14).Value = -4219
$ws.Cells.Item(48, 8).Value = 21747.5
$ws.Cells.Item(48, 9).Value = 6000
$ws.Cells.Item(48, 10).Value = 37495
$ws.Cells.Item(48, 11).Value = 6000
$ws.Cells.Item(48, 12).Value = 37495
$ws.Cells.Item(48, 13).Value = -5339
$ws.Cells.Item(48, 14).Value = -38817
$ws.Cells.Item(55, 8).Value = 554.9286
$ws.Cells.Item(55, 9).Value = 626.44446
$ws.Cells.Item(55, 10).Value = 426.2
$ws.Cells.Item(55, 11).Value = 626.44446
$ws.Cells.Item(55, 12).Value = 426.2
$ws.Cells.Item(55, 13).Value = -453.44446
$ws.Cells.Item(55, 14).Value = -772.2
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).ClearContents()
$ws.Cells.Item(68, 14).ClearContents()
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).ClearContents()
$ws.Cells.Item(71, 14).ClearContents()
$ws.Cells.Item(104, 8).Value = 27370
$ws.Cells.Item(104, 10).Value = 27370
$ws.Cells.Item(104, 12).Value = 27370
$ws.Cells.Item(104, 14).Value = -34358
$ws.Cells.Item(122, 8).Value = 2844195.8
$ws.Cells.Item(122, 9).Value = 3273.394
$ws.Cells.Item(122, 10).Value = 11366963
$ws.Cells.Item(122, 11).Value = 9820.181999999999
$ws.Cells.Item(122, 12).Value = 34100889
$ws.Cells.Item(122, 13).Value = -7370.181999999999
$ws.Cells.Item(122, 14).Value = -34105789
$ws.Cells.Item(132, 8).Value = 1630.7241
$ws.Cells.Item(132, 9).Value = 1124.6
$ws.Cells.Item(132, 10).Value = 2755.4443
$ws.Cells.Item(132, 11).Value = 3373.8
$ws.Cells.Item(132, 12).Value = 8266.332900000001
$ws.Cells.Item(132, 13).Value = -843.7999999999997
$ws.Cells.Item(132, 14).Value = -13326.3329

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(62, 8).Value = 7610.4443
$ws.Cells.Item(62, 9).Value = 5918.8
$ws.Cells.Item(62, 10).Value = 9725
$ws.Cells.Item(62, 11).Value = 5918.8
$ws.Cells.Item(62, 12).Value = 9725
$ws.Cells.Item(62, 13).Value = -5294.8
$ws.Cells.Item(62, 14).Value = -10973
$ws.Cells.Item(65, 8).Value = 7610.4443
$ws.Cells.Item(65, 9).Value = 5918.8
$ws.Cells.Item(65, 10).Value = 9725
$ws.Cells.Item(65, 11).Value = 29594
$ws.Cells.Item(65, 12).Value = 48625
$ws.Cells.Item(65, 13).Value = -26474
$ws.Cells.Item(65, 14).Value = -54865
$ws.Cells.Item(81, 8).Value = 4944.205
$ws.Cells.Item(81, 9).Value = 2992.8262
$ws.Cells.Item(81, 11).Value = 5985.6524
$ws.Cells.Item(81, 13).Value = -4924.6524
$ws.Cells.Item(84, 8).Value = 4944.205
$ws.Cells.Item(84, 9).Value = 2992.8262
$ws.Cells.Item(84, 11).Value = 29928.262
$ws.Cells.Item(84, 13).Value = -24624.262
$ws.Cells.Item(122, 8).Value = 14290720
$ws.Cells.Item(122, 9).Value = 4058.3
$ws.Cells.Item(122, 10).Value = 50007376
$ws.Cells.Item(122, 11).Value = 12174.9
$ws.Cells.Item(122, 12).Value = 150022128
$ws.Cells.Item(122, 13).Value = -9724.900000000001
$ws.Cells.Item(122, 14).Value = -150027028
$ws.Cells.Item(126, 8).Value = 11908464
$ws.Cells.Item(126, 9).Value = 11908464
$ws.Cells.Item(126, 11).Value = 35725392
$ws.Cells.Item(126, 13).Value = -35722922
$ws.Cells.Item(132, 8).Value = 2492.44
$ws.Cells.Item(132, 9).Value = 2309.6135
$ws.Cells.Item(132, 11).Value = 6928.8405
$ws.Cells.Item(132, 13).Value = -4398.8405
$ws.Cells.Item(136, 8).Value = 4497.647
$ws.Cells.Item(136, 9).Value = 3035.5386
$ws.Cells.Item(136, 10).Value = 9249.5
$ws.Cells.Item(136, 11).Value = 9106.6158
$ws.Cells.Item(136, 12).Value = 27748.5
$ws.Cells.Item(136, 13).Value = -6556.6158
$ws.Cells.Item(136, 14).Value = -32848.5
